# Scheduled market-data refresh: update currentAveragePrice and the
# derived profit columns (H:N) across the leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3964
$ws.Range("I40").Value = 1974.75
$ws.Range("J40").Value = 4282.28
$ws.Range("K40").Value = 1974.75
$ws.Range("L40").Value = 4282.28
$ws.Range("M40").Value = -1799.75
$ws.Range("N40").Value = -4632.28

$ws.Range("H58").Value = 564.5
$ws.Range("I58").Value = 86.166664
$ws.Range("J58").Value = 1999.5
$ws.Range("K58").Value = 258.499992
$ws.Range("L58").Value = 5998.5
$ws.Range("M58").Value = -108.499992
$ws.Range("N58").Value = -6298.5

$ws.Range("H62").Value = 7570.3335
$ws.Range("I62").Value = 996.5
$ws.Range("K62").Value = 996.5
$ws.Range("M62").Value = -372.5

$ws.Range("H65").Value = 7570.3335
$ws.Range("I65").Value = 996.5
$ws.Range("K65").Value = 4982.5
$ws.Range("M65").Value = -1862.5

$ws.Range("H135").Value = 1691.5938
$ws.Range("I135").Value = 778.4091
$ws.Range("K135").Value = 7005.6819
$ws.Range("M135").Value = -4470.6819

$ws.Range("H137").Value = 3517.7666
$ws.Range("I137").Value = 2482.3845
$ws.Range("J137").Value = 4309.5293
$ws.Range("K137").Value = 7447.1535
$ws.Range("L137").Value = 12928.5879
$ws.Range("M137").Value = -4897.1535
$ws.Range("N137").Value = -18028.5879

$ws.Range("H138").Value = 2353.87
$ws.Range("I138").Value = 1417.5807
$ws.Range("K138").Value = 4252.742099999999
$ws.Range("M138").Value = 887.2579000000005

$ws.Range("H141").Value = 1721.7142
$ws.Range("I141").Value = 1388.6471
$ws.Range("K141").Value = 4165.9413
$ws.Range("M141").Value = 1014.0587

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4416.548
$ws.Range("I32").Value = 2829.491
$ws.Range("J32").Value = 9265.888999999999
$ws.Range("K32").Value = 2829.491
$ws.Range("L32").Value = 9265.888999999999
$ws.Range("M32").Value = -2542.491
$ws.Range("N32").Value = -9839.888999999999

$ws.Range("H122").Value = 2997.5862
$ws.Range("I122").Value = 2387.3809
$ws.Range("K122").Value = 7162.1427
$ws.Range("M122").Value = -4712.1427

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2938.8572
$ws.Range("I20").Value = 2457.6365
$ws.Range("K20").Value = 2457.6365
$ws.Range("M20").Value = -2210.6365

$ws.Range("H99").Value = 87011.75
$ws.Range("I99").Value = 252302
$ws.Range("J99").Value = 4366.625
$ws.Range("K99").Value = 252302
$ws.Range("L99").Value = 4366.625
$ws.Range("M99").Value = -250804
$ws.Range("N99").Value = -7362.625

$ws.Range("H134").Value = 2546.359
$ws.Range("I134").Value = 798.59375
$ws.Range("K134").Value = 2395.78125
$ws.Range("M134").Value = 139.21875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3027.6765
$ws.Range("I31").Value = 2064.6667
$ws.Range("K31").Value = 2064.6667
$ws.Range("M31").Value = -1769.6667

$ws.Range("H34").Value = 3027.6765
$ws.Range("I34").Value = 2064.6667
$ws.Range("K34").Value = 2064.6667
$ws.Range("M34").Value = -1862.6667

$ws.Range("H59").Value = 39999.8
$ws.Range("J59").Value = 39999.8
$ws.Range("L59").Value = 39999.8
$ws.Range("N59").Value = -42289.8

$ws.Range("H122").Value = 3429.5833
$ws.Range("I122").Value = 3286
$ws.Range("K122").Value = 9858
$ws.Range("M122").Value = -7408

$ws.Range("H134").Value = 37682.617
$ws.Range("I134").Value = 55285.824
$ws.Range("J134").Value = 4432.1113
$ws.Range("K134").Value = 165857.472
$ws.Range("L134").Value = 13296.3339
$ws.Range("M134").Value = -163322.472
$ws.Range("N134").Value = -18366.3339

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 4641.364
$ws.Range("I33").Value = 82.2
$ws.Range("K33").Value = 493.2
$ws.Range("M33").Value = -210.2

$ws.Range("H86").Value = 430
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 430
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H141").Value = 2849.8
$ws.Range("I141").Value = 1701.6
$ws.Range("J141").Value = 3998
$ws.Range("K141").Value = 5104.799999999999
$ws.Range("L141").Value = 11994
$ws.Range("M141").Value = 75.20000000000073
$ws.Range("N141").Value = -22354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 647.04346
$ws.Range("I107").Value = 669.8823
$ws.Range("J107").Value = 582.3333
$ws.Range("K107").Value = 669.8823
$ws.Range("L107").Value = 582.3333
$ws.Range("M107").Value = 1250.1177
$ws.Range("N107").Value = -4422.3333

$ws.Range("H113").Value = 2703.6
$ws.Range("I113").Value = 1301.2
$ws.Range("K113").Value = 1301.2
$ws.Range("M113").Value = 868.8

$ws.Range("H122").Value = 5320.1
$ws.Range("I122").Value = 4348.4287
$ws.Range("K122").Value = 13045.2861
$ws.Range("M122").Value = -10595.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 20000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 20000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 20000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -20590

$ws.Range("H46").Value = 8011.375
$ws.Range("I46").Value = 7209.3335
$ws.Range("J46").Value = 9042.571
$ws.Range("K46").Value = 7209.3335
$ws.Range("L46").Value = 9042.571
$ws.Range("M46").Value = -7021.3335
$ws.Range("N46").Value = -9418.571

$ws.Range("H100").Value = 3242.5806
$ws.Range("I100").Value = 3059.0908
$ws.Range("J100").Value = 3691.111
$ws.Range("K100").Value = 3059.0908
$ws.Range("L100").Value = 3691.111
$ws.Range("M100").Value = -2518.0908
$ws.Range("N100").Value = -4773.111

$ws.Range("H122").Value = 6715.6
$ws.Range("I122").Value = 3794.8
$ws.Range("K122").Value = 11384.4
$ws.Range("M122").Value = -8934.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 9571.429
$ws.Range("J52").Value = 9750
$ws.Range("L52").Value = 9750
$ws.Range("N52").Value = -10202

$ws.Range("H69").Value = 17237.5
$ws.Range("J69").Value = 17237.5
$ws.Range("L69").Value = 17237.5
$ws.Range("N69").Value = -18735.5

$ws.Range("H72").Value = 17237.5
$ws.Range("J72").Value = 17237.5
$ws.Range("L72").Value = 51712.5
$ws.Range("N72").Value = -59200.5

$ws.Range("H132").Value = 5337.343
$ws.Range("I132").Value = 5591.357
$ws.Range("J132").Value = 4321.2856
$ws.Range("K132").Value = 16774.071
$ws.Range("L132").Value = 12963.8568
$ws.Range("M132").Value = -14244.071
$ws.Range("N132").Value = -18023.8568

$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360
